$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the trailing-space typo in the "correct_answer" column for question 8 (row 9):
# "c " -> "c"
$ws.Range("F9").Value = "c"

# Update the active selection on the sheet to match the saved view state
$ws.Range("G9").Select()
